$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.457.24'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '1.751.06'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'322.49"
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = "'0.4250"
$ws.Range('E7').Value = '  -4.60%  '
$ws.Range('D8').Value = "'0.3598"
$ws.Range('E8').Value = '  -3.38%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'42.30"
$ws.Range('E9').Value = '  -5.97%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07461"
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('D11').Value = "'1.098"
$ws.Range('E11').Value = '  -3.39%  '
$ws.Range('D12').Value = "'1.000"
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E13').Value = '  -6.96%  '
$ws.Range('D14').Value = "'6.031"
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').Value = "'7.213"
$ws.Range('E15').Value = '  -5.45%  '
$ws.Range('D16').Value = '1.746.24'
$ws.Range('E16').Value = '  -5.34%  '
$ws.Range('D17').Value = "'93.13"
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('E18').Value = '  -2.28%  '
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').Value = "'0.9992"
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  -3.13%  '
$ws.Range('D22').Value = "'5.890"
$ws.Range('E22').Value = '  -5.92%  '
$ws.Range('D23').Value = '27.507.10'
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('D24').Value = "'11.23"
$ws.Range('E24').Value = '  -4.07%  '
$ws.Range('E25').Value = '  -5.25%  '
$ws.Range('D26').Value = "'162.06"
$ws.Range('E26').Value = '  +3.35%  '
$ws.Range('D27').Value = "'20.18"
$ws.Range('E27').Value = '  -3.36%  '
$ws.Range('D28').Value = '1.945.01'
$ws.Range('E28').Value = '  -5.17%  '
$ws.Range('D29').Value = "'2.138"
$ws.Range('E29').Value = '  -8.18%  '
$ws.Range('D30').Value = "'123.87"
$ws.Range('E30').Value = '  -4.07%  '
$ws.Range('D31').Value = "'1.105"
$ws.Range('E31').Value = '  -8.96%  '
$ws.Range('D32').Value = "'3.661"
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').Value = "'5.552"
$ws.Range('E33').Value = '  -6.98%  '
$ws.Range('D34').Value = "'0.08899"
$ws.Range('E34').Value = '  -3.91%  '
$ws.Range('E35').Value = '  -7.72%  '
$ws.Range('D36').Value = "'0.02290"
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('D37').Value = "'0.2090"
$ws.Range('E37').Value = '  -4.98%  '
$ws.Range('E38').Value = '  -4.14%  '
$ws.Range('E39').Value = '  -4.82%  '
$ws.Range('D40').Value = "'4.955"
$ws.Range('E40').Value = '  -5.13%  '
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('D42').Value = "'7.929"
$ws.Range('D43').Value = "'0.9997"
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = "'1.390"
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('D45').Value = "'13.32"
$ws.Range('E45').Value = '  -5.00%  '
$ws.Range('D46').Value = "'0.5877"
$ws.Range('E46').Value = '  -4.85%  '
$ws.Range('D47').Value = "'3.689"
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').Value = "'123.04"
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('D49').Value = "'1.964"
$ws.Range('E49').Value = '  -4.28%  '
$ws.Range('D50').Value = "'1.165"
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = "'0.06828"
$ws.Range('E51').Value = '  -2.63%  '
